$d = $word.ActiveDocument

# Locate the paragraph that ends the payment-recording bullet, which is
# where the three new "Get analytical report" bullets are appended after.
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i).Range
    if ($cand.Text -like "*Record the payment methods and amounts*") {
        $idx = $i
        break
    }
}

$paymentPara = $d.Paragraphs($idx).Range

# Insert the first new bullet ("Get analytical report") at level 2 (ilvl=1).
$paymentPara.InsertParagraphAfter()
$p1 = $d.Paragraphs($idx + 1).Range
$p1.Text = "Get analytical report"
$p1.ListFormat.ListLevelNumber = 2

# Insert the second new bullet at level 3 (ilvl=2).
$p1.InsertParagraphAfter()
$p2 = $d.Paragraphs($idx + 2).Range
$p2.Text = "Manager requests analytics from database"
$p2.ListFormat.ListLevelNumber = 3

# Insert the third new bullet at level 3 (ilvl=2).
$p2.InsertParagraphAfter()
$p3 = $d.Paragraphs($idx + 3).Range
$p3.Text = "Relevant data is queried and supplied to manager"
$p3.ListFormat.ListLevelNumber = 3

# Move the "_GoBack" bookmark (Word keeps only one) from wherever the last
# edit left it onto the point between "G" and "et analytical report" in the
# freshly-inserted first bullet. Adding it here automatically removes it
# from its old location, matching Word's single-_GoBack behavior.
$p1reload = $d.Paragraphs($idx + 1).Range
$bmRange = $d.Range($p1reload.Start + 1, $p1reload.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
